# Auto-generated Excel COM-interop script applying the diff changes
$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item('展览')
$ws.Range("F2").Value = 1543
$ws.Range("F3").Value = 905
$ws.Range("F4").Value = 467
$ws.Range("F5").Value = 912
$ws.Range("C7").Value = '杭州·AD04动漫展'
$ws.Range("E7").Value = '2024.07.13 10:00-07.14 17:00'
$ws.Range("F7").Value = 7805
$ws.Range("G7").Value = 75
$ws.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=85012'
$ws.Range("I7").Value = '//i0.hdslb.com/bfs/openplatform/202405/y1iKqqnh1715326769523.jpeg'
$ws.Range("C8").Value = '杭州·AD04动漫展-Pile·内场票'
$ws.Range("F8").Value = 139
$ws.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=86819'
$ws.Range("I8").Value = '//i0.hdslb.com/bfs/openplatform/202406/AULfssPJ1717482529866.jpeg'
$ws.Range("C9").Value = '杭州·AD04动漫展·小泽亚李·内场票'
$ws.Range("D9").Value = '阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心'
$ws.Range("E9").Value = '2024.07.13 09:30-07.13 17:00'
$ws.Range("F9").Value = 36
$ws.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=87241'
$ws.Range("I9").Value = '//i0.hdslb.com/bfs/openplatform/202406/PtX41aGD1718245480604.jpeg'
$ws.Range("C10").Value = '杭州·代号鸢only-广陵大学'
$ws.Range("D10").Value = '康候圣街99号 顺丰创新中心'
$ws.Range("E10").Value = '2024.07.13 09:30-07.13 17:30'
$ws.Range("F10").Value = 1937
$ws.Range("G10").Value = '已售罄'
$ws.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=83289'
$ws.Range("I10").Value = '//i0.hdslb.com/bfs/openplatform/202406/k2CTTf491718604574410.jpeg'
$ws.Range("C11").Value = '杭州·草莓动漫节'
$ws.Range("E11").Value = '2024.07.13 09:00-07.14 17:00'
$ws.Range("F11").Value = 5641
$ws.Range("G11").Value = 70
$ws.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=84229'
$ws.Range("I11").Value = '//i1.hdslb.com/bfs/openplatform/202406/czPRn1ve1718875288240.jpeg'
$ws.Range("B12").Value = "'2024-07-13"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = '杭州·草莓动漫节内场票·钱文青'
$ws.Range("D12").Value = '松合路2号 钱塘文体中心'
$ws.Range("E12").Value = '2024.07.13 09:00-07.13 17:00'
$ws.Range("F12").Value = 575
$ws.Range("G12").Value = 238
$ws.Range("H12").Value = 'https://show.bilibili.com/platform/detail.html?id=84851'
$ws.Range("I12").Value = '//i2.hdslb.com/bfs/openplatform/202406/Vx16dKjM1718875179041.jpeg'
$ws.Range("F14").Value = 7876
$ws.Range("G14").Value = 75
$ws.Range("F15").Value = 9219
$ws.Range("F16").Value = 1154
$ws.Range("F17").Value = 917
$ws.Range("F18").Value = 4517
$ws.Range("F19").Value = 682
$ws.Range("F20").Value = 256
$ws.Range("F22").Value = 291
$ws.Range("F24").Value = 1204
$ws.Range("F25").Value = 124
$ws.Range("F26").Value = 1694
$ws.Range("F27").Value = 734
$ws.Range("F28").Value = 955
$ws.Range("F29").Value = 13
$ws.Range("F30").Value = 1895
$ws.Range("F31").Value = 345
$ws.Range("F32").Value = 2333
$ws.Range("F34").Value = 1490
$ws.Range("F35").Value = 72
$ws.Range("F37").Value = 4
$ws.Range("F38").Value = 805
$ws.Range("F39").Value = 517
$ws.Range("F40").Value = 3005
$ws.Range("F41").Value = 4150
$ws.Range("F42").Value = 194
$ws.Range("F43").Value = 51
$ws.Range("F44").Value = 430
$ws.Range("F47").Value = 862
$ws.Range("F48").Value = 179
$ws.Range("F49").Value = 4101

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item('演出')
$ws.Range("F15").Value = 49

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item('本地生活')
$ws.Range("F2").Value = 5324

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item('全部类型')
$ws.Range("F3").Value = 1543
$ws.Range("F4").Value = 905
$ws.Range("F5").Value = 467
$ws.Range("F6").Value = 912
$ws.Range("C8").Value = '杭州·【早鸟6折】《忱宴·渐渐被你吸引》热血动漫二次元ACG演唱会'
$ws.Range("D8").Value = '湖墅南路136-138号 浙话艺术剧院'
$ws.Range("E8").Value = '2024.07.13 19:30-07.13 21:30'
$ws.Range("F8").Value = 24
$ws.Range("G8").Value = 60
$ws.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=85011'
$ws.Range("I8").Value = '//i1.hdslb.com/bfs/openplatform/202404/2Gd8eLva1714379746993.jpeg'
$ws.Range("C9").Value = '杭州·海上钢琴师—一生必听的电影名曲《泰坦尼克号》《花样年华》《海上钢琴师》'
$ws.Range("D9").Value = '曙光路31号 浙江音乐厅'
$ws.Range("E9").Value = '2024.07.13 19:30-07.13 21:00'
$ws.Range("F9").Value = 8
$ws.Range("G9").Value = 100
$ws.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=85889'
$ws.Range("I9").Value = '//i0.hdslb.com/bfs/openplatform/202405/52kxbBTh1716096935602.jpeg'
$ws.Range("C10").Value = '杭州·草莓动漫节'
$ws.Range("D10").Value = '松合路2号 钱塘文体中心'
$ws.Range("E10").Value = '2024.07.13 09:00-07.14 17:00'
$ws.Range("F10").Value = 5641
$ws.Range("G10").Value = 70
$ws.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=84229'
$ws.Range("I10").Value = '//i1.hdslb.com/bfs/openplatform/202406/czPRn1ve1718875288240.jpeg'
$ws.Range("C11").Value = '杭州·草莓动漫节内场票·钱文青'
$ws.Range("E11").Value = '2024.07.13 09:00-07.13 17:00'
$ws.Range("F11").Value = 575
$ws.Range("G11").Value = 238
$ws.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=84851'
$ws.Range("I11").Value = '//i2.hdslb.com/bfs/openplatform/202406/Vx16dKjM1718875179041.jpeg'
$ws.Range("B12").Value = "'2024-07-20"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = '【会员购严选】杭州·首届次元格子动漫展-进入格子空间，探索次元世界！'
$ws.Range("D12").Value = '钱江世纪城奔竞大道353号 杭州国际博览中心'
$ws.Range("E12").Value = '2024.07.20 09:00-07.22 17:00'
$ws.Range("F12").Value = 7877
$ws.Range("G12").Value = 75
$ws.Range("H12").Value = 'https://show.bilibili.com/platform/detail.html?id=85616'
$ws.Range("I12").Value = '//i1.hdslb.com/bfs/openplatform/202405/5Dne5VqI1715753018080.jpeg'
$ws.Range("C13").Value = '杭州·TCD国潮动漫游戏嘉年华'
$ws.Range("D13").Value = '阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心'
$ws.Range("E13").Value = '2024.07.20 09:30-07.21 17:00'
$ws.Range("F13").Value = 9219
$ws.Range("G13").Value = 65
$ws.Range("H13").Value = 'https://show.bilibili.com/platform/detail.html?id=85699'
$ws.Range("I13").Value = '//i1.hdslb.com/bfs/openplatform/202406/QzaksReK1718190369702.jpeg'
$ws.Range("F15").Value = 917
$ws.Range("F16").Value = 4518
$ws.Range("F17").Value = 682
$ws.Range("F18").Value = 256
$ws.Range("F20").Value = 291
$ws.Range("F23").Value = 1204
$ws.Range("F24").Value = 124
$ws.Range("F25").Value = 1694
$ws.Range("F26").Value = 734
$ws.Range("F27").Value = 955
$ws.Range("F28").Value = 13
$ws.Range("F29").Value = 1895
$ws.Range("F30").Value = 345
$ws.Range("F31").Value = 2333
$ws.Range("F33").Value = 72
$ws.Range("F36").Value = 805
$ws.Range("F39").Value = 517
$ws.Range("F40").Value = 4150
$ws.Range("F42").Value = 194
$ws.Range("F43").Value = 51
$ws.Range("F44").Value = 430
$ws.Range("F47").Value = 862
$ws.Range("F48").Value = 179
$ws.Range("F49").Value = 4101
